$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (shifts J..R -> K..S)
$ws.Columns("J").Insert()

# New header for the inserted column
$ws.Cells.Item(1, 10).Value = "netto_da_pagare"

# New data values for netto_da_pagare, rows 2-94 (column J / index 10)
$nettoValues = @(255,1339,1695,1342,755,402,1271,1306,1264,1262,2472,1259,1306,1262,1260,1963,1960,2104,1193,1861,2210,2130,3158,2141,2155,2074,2083,1928,1976,2111,2465,1326,2300,2412,2123,2108,2119,2143,2145,2144,2163,2414,2276,2530,1362,1970,2189,2793,2191,2194,2261,2837,2891,2655,2586,2763,2367,1383,2755,1895,2444,2375,2242,2209,2273,2001,1999,2001,1871,1838,1818,2050,2058,1977,1980,2071,3960,1929,2007,2193,2174,2056,1852,1846,2090,2128,2078,2141,2078,1913,2321,2087,2302)

for ($i = 0; $i -lt $nettoValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $nettoValues[$i]
}
